# "Update displayed press brands"
# The sheet maps press titles (col A) to a press "Typ" (col B). Four titles
# are being re-flagged as "NIEUWZGLĘDNIONE" (disregarded), reusing the
# shared string already present in the workbook (B20 already has this value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B15").Value  = "NIEUWZGLĘDNIONE"
$ws.Range("B58").Value  = "NIEUWZGLĘDNIONE"
$ws.Range("B91").Value  = "NIEUWZGLĘDNIONE"
$ws.Range("B104").Value = "NIEUWZGLĘDNIONE"

# Bring the view to where the author was working (scrolled further down,
# with a couple of cells picked out) and re-select accordingly.
$ws.Activate()

$win = $excel.ActiveWindow
$win.ScrollRow    = 94
$win.ScrollColumn = 1

$ws.Range("E107,O100").Select()
$ws.Range("O100").Activate()
